$wb = $excel.ActiveWorkbook

# Select A2 on "Ville Accroche" before switching away (matches final selection state)
$wsVilleAccroche = $wb.Worksheets.Item("Ville Accroche")
$wsVilleAccroche.Activate()
$wsVilleAccroche.Range("A2").Select()

# Create the new first sheet "Département Accroche"
$new = $wb.Worksheets.Add()
$new.Name = "Département Accroche"

# Row 1 header
$new.Range("A1").Value = "phrase1"
$new.Range("B1").Value = "phrase2"
$new.Range("C1").Value = "phrase3"

# row 2
$new.Range("B2").Value = "Grâce à notre site, vous pouvez trouver rapidement"
$new.Range("C2").Value = "dès que vous en ressentez la nécessité"
$new.Range("A2").Value = "Vous vous trouvez sur la page de notre annuaire"

# row 3
$new.Range("B3").Value = "Découvrez les solutions grâce à nos conseillers afin de trouver"
$new.Range("C3").Value = "dès que le besoin se fait ressentir"
$new.Range("A3").Value = "Bienvenu sur la page de l'annuaire"

# row 4
$new.Range("B4").Value = "Toutes nos équipes vous accompagnent afin de trouver"
$new.Range("C4").Value = "à tout moment de la journée"
$new.Range("A4").Value = "Vous êtes bien sur l'annuaire"

# row 5
$new.Range("B5").Value = "Tous nos consultants spécialisées vous aident à trouver"
$new.Range("C5").Value = "lorsque vous avez besoin d'une consultation"
$new.Range("A5").Value = "L'accès à l'annuaire de cette page web"

# row 6
$new.Range("B6").Value = "Tous nos conseillers vous accompagnent afin de récupérer et trouver"
$new.Range("C6").Value = "si vous avez un besoin médical en dehors des horaires habituels"
$new.Range("A6").Value = "L'annuaire de cette page vous donne accès à"

# row 7
$new.Range("B7").Value = "Grâce à nos équipes de téléconseillers, trouvez facilement"
$new.Range("C7").Value = "dès que le besoin demande un conseil médical"
$new.Range("A7").Value = "Bonjour, nous sommes ravis de vous accueillir sur l'annuaire"

# row 8
$new.Range("B8").Value = "Nos consultants en ligne, vous aident pour trouver rapidement "
$new.Range("C8").Value = "dès que vous ressentez l'envie de parler à un spécialiste"
$new.Range("A8").Value = "Vous êtes bien sur notre annuaire"

# row 9
$new.Range("B9").Value = "Une question médicale, besoin d'un conseil nos équipes vous trouvent"
$new.Range("C9").Value = "dés que le moment vous parait opportun"
$new.Range("A9").Value = "Bienvenu sur l'annuaire"

# row 10
$new.Range("B10").Value = "Un besoin urgent, une demande médicale, nous vous aidons à trouver"
$new.Range("C10").Value = "lorsque vous pensez avoir besoin d'un conseil médical"
$new.Range("A10").Value = "Bienvenu sur notre page annuaire"

# row 11
$new.Range("B11").Value = "Notre offre de renseignements, vous permet facilement de trouver"
$new.Range("C11").Value = "si vous pensez que c'est le moment d'avoir un conseil d'un spécialiste"
$new.Range("A11").Value = "Bienvenu sur la page de notre annuaire"

# row 12
$new.Range("B12").Value = "Notre offre de renseignements, vous permet rapidement de récupérer"
$new.Range("C12").Value = "dès qu'un besoin urgent se fait ressentir"
$new.Range("A12").Value = "Bienvenu sur notre annuaire"

# row 13
$new.Range("B13").Value = "Avec notre site internet, vous pouvez récupérer facilement"
$new.Range("C13").Value = "dès lors que le besoin d'un conseil médical se fait ressentir"
$new.Range("A13").Value = "Bienvenu sur votre annuaire"

# row 14
$new.Range("B14").Value = "Découvrez les différentes solutions afin de trouver rapidement"
$new.Range("C14").Value = "dès que vous en éprouvez le besoin"
$new.Range("A14").Value = "Vous êtes bien arrivés sur l'annuaire"

# row 15
$new.Range("B15").Value = "Nos équipes vous accompagnent dans la recherche afin de trouver"
$new.Range("C15").Value = "qui pourra répondre à vos demandes médicales"
$new.Range("A15").Value = "Vous êtes bien sur la page annuaire"

# row 16
$new.Range("B16").Value = "Avec notre site internet, vous pouvez récupérer facilement"
$new.Range("C16").Value = "qui pourra vous aider dans vos demandes médicales"
$new.Range("A16").Value = "Vous êtes bien sur la page de l'annuaire "

# row 17
$new.Range("B17").Value = "En vous connectant à notre page, vous trouverez facilement comment récupérer"
$new.Range("C17").Value = "qui pourra vous rassurer sur vos problèmes d'ordre médicaux"
$new.Range("A17").Value = "Bonjour, vous voilà bien sur l'annuaire"

# row 18
$new.Range("B18").Value = "Toutes nos équipes spécialisées vous aident afin de trouver"
$new.Range("C18").Value = "qui pourra vous rassurer sur vos symptômes"
$new.Range("A18").Value = "Bonjour, vous êtes arrivés sur l'annuaire"

# row 19
$new.Range("B19").Value = "Grâce à nos téléconseillers, découvrez comment trouver "
$new.Range("C19").Value = "qui pourra analyser vos symptômes"
$new.Range("A19").Value = "Bonjour et bienvenu sur notre annuaire "

# row 20
$new.Range("B20").Value = "Grâce à nos équipes de téléconseillers, découvrez comment récupérer"
$new.Range("C20").Value = "qui vous aidera afin d'analyser vos symptômes"
$new.Range("A20").Value = "Bonjour, vous êtes bien sur l'annuaire"

# row 21
$new.Range("B21").Value = "Grâce à nos consultants téléconseillers, découvrez comment chercher et trouver "
$new.Range("C21").Value = "à tout moment de la semaine"
$new.Range("A21").Value = "Vous êtes à présent sur la page de notre annuaire"

# row 22
$new.Range("B22").Value = "Cette page web vous offre un service qui vous permettra facilement de trouver"
$new.Range("C22").Value = "qui sera à même de répondre à toutes vos questions"

# row 23
$new.Range("B23").Value = "Cette page internet vous donne accès à un service afin de récupérer "
$new.Range("C23").Value = "qui sera capable de répondre à vos questions et analyser vos symptômes"

# row 24
$new.Range("B24").Value = "Cette page web vous donne un service de mise en relation afin de trouver"
$new.Range("C24").Value = "qui sera dans la capacité de vous aider et répondre à vos questions médicales"

# row 25
$new.Range("B25").Value = "En passant par notre site internet, nous vous aidons rapidement afin de trouver"
$new.Range("C25").Value = "qui vous répondra sur les questions et symptômes que vous ressentez"

# row 26
$new.Range("C26").Value = "qui écoutera les symptômes que vous ressentez"

# row 27
$new.Range("C27").Value = "qui pourra vous donner des élements de réponses sur vos symptômes"

# row 28
$new.Range("C28").Value = "qui vous dira ce qu'il faut faire selon votre situation médicale"

# row 29
$new.Range("C29").Value = "qui vous écoutera afin de trouver une solution à vos problèmes"

# row 30
$new.Range("C30").Value = "qui donnera des élements de réponses à vos questions"

# row 31
$new.Range("C31").Value = "qui sera à l'écoute de vos demandes afin de trouver une solution"

# Move the new sheet to the first position and make it active
$new.Move($wb.Worksheets.Item(1))
$first = $wb.Worksheets.Item(1)
$first.Activate()
$first.Range("A25").Select()